$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results (loading_percent) for the 380 kV case.
# Columns B, C, E, F, G, J, L, N are updated for each row (A holds the row index 0..23).
$data = @{
    2 = @{ "B"=28.38421571448552; "C"=11.17112742013186; "E"=9.966375593026774; "F"=50.95377225451496; "G"=3.75584773364357; "J"=10.53001651883042; "L"=11.60666337330299; "N"=21.58009586826574 }
    3 = @{ "B"=28.01066546359994; "C"=10.66285483182122; "E"=9.945332846457788; "F"=50.78604654158608; "G"=3.760749278433119; "J"=10.55943188511145; "L"=11.61670669224455; "N"=21.63668320938842 }
    4 = @{ "B"=27.78770282930379; "C"=10.34204414884653; "E"=9.932240218206308; "F"=50.69555784195282; "G"=3.763911430104563; "J"=10.57864476482488; "L"=11.62471092831302; "N"=21.67342288326188 }
    5 = @{ "B"=27.69856094908152; "C"=10.20932920110105; "E"=9.926860680762474; "F"=50.66183917060201; "G"=3.76523856552908; "J"=10.58676387909527; "L"=11.62843527971425; "N"=21.68889595666347 }
    6 = @{ "B"=27.68386580757438; "C"=10.18717905162543; "E"=9.925964718891962; "F"=50.65643117572078; "G"=3.765461267476163; "J"=10.58812955278169; "L"=11.62908165511815; "N"=21.69149552255223 }
    7 = @{ "B"=27.78649353770869; "C"=10.34026203164988; "E"=9.932167847927381; "F"=50.69509030487953; "G"=3.763929172078905; "J"=10.57875308886067; "L"=11.62475928275; "N"=21.67362952885555 }
    8 = @{ "B"=28.25415356007875; "C"=10.9978061397709; "E"=9.959155292723725; "F"=50.89335654416418; "G"=3.757506216979974; "J"=10.53991995324069; "L"=11.60974523088471; "N"=21.59919300205521 }
    9 = @{ "B"=29.21708941699994; "C"=12.20998092431624; "E"=10.01075084981769; "F"=51.38051970308705; "G"=3.746113986542019; "J"=10.47290246704583; "L"=11.59485954363822; "N"=21.46906223342123 }
    10 = @{ "B"=29.94578383048765; "C"=13.04446631763255; "E"=10.04790747198954; "F"=51.79712355238038; "G"=3.738467067137738; "J"=10.42922764711887; "L"=11.5927616285874; "N"=21.38312584527427 }
    11 = @{ "B"=30.28045088870332; "C"=13.41044953677648; "E"=10.06466066263286; "F"=51.99904356565127; "G"=3.735143008652261; "J"=10.41056600486396; "L"=11.5937168388867; "N"=21.34613496303335 }
    12 = @{ "B"=30.40751822963791; "C"=13.5469807167235; "E"=10.07098412766409; "F"=52.07725438637305; "G"=3.733906324977831; "J"=10.40367279894361; "L"=11.59435211843758; "N"=21.33243023566319 }
    13 = @{ "B"=30.38013895852972; "C"=13.51766925415571; "E"=10.06962315837254; "F"=52.06033314184082; "G"=3.734171688142556; "J"=10.40514965358923; "L"=11.59420315047503; "N"=21.33536831073034 }
    14 = @{ "B"=30.29089874407056; "C"=13.42172380936643; "E"=10.06518130567946; "F"=52.00544319832371; "G"=3.735040824645712; "J"=10.4099954177817; "L"=11.5937636274253; "N"=21.34500139208298 }
    15 = @{ "B"=30.2362767573507; "C"=13.36268368298853; "E"=10.06245788610794; "F"=51.97204806250786; "G"=3.735576064606999; "J"=10.41298619319256; "L"=11.59353000023638; "N"=21.35094140132613 }
    16 = @{ "B"=29.92396614823252; "C"=13.02026538413177; "E"=10.04680967917023; "F"=51.78417424250482; "G"=3.738687398632701; "J"=10.43047150985424; "L"=11.59273754308386; "N"=21.38558563882781 }
    17 = @{ "B"=29.73310149162702; "C"=12.80663664818326; "E"=10.03717234796867; "F"=51.67207253441136; "G"=3.74063557417365; "J"=10.44150719625468; "L"=11.5927397832701; "N"=21.40737758045353 }
    18 = @{ "B"=29.62362769146716; "C"=12.6824850916694; "E"=10.03161508017743; "F"=51.60876454631519; "G"=3.741770671083263; "J"=10.44796816234116; "L"=11.59292087760686; "N"=21.42010953865235 }
    19 = @{ "B"=29.5866181073846; "C"=12.64023330108089; "E"=10.02973103069325; "F"=51.58753152427376; "G"=3.742157500622206; "J"=10.45017523095427; "L"=11.59301310442018; "N"=21.42445431579455 }
    20 = @{ "B"=29.7533885220852; "C"=12.82951082379363; "E"=10.03819971420294; "F"=51.68388511581702; "G"=3.740426682044055; "J"=10.44032067793946; "L"=11.59272094269062; "N"=21.40503731395048 }
    21 = @{ "B"=30.31710258929564; "C"=13.44996189534287; "E"=10.06648653750366; "F"=52.02151856274769; "G"=3.734784940626865; "J"=10.40856738832391; "L"=11.59388531037055; "N"=21.34216369608973 }
    22 = @{ "B"=30.68742732887587; "C"=13.84342763721294; "E"=10.084853933265; "F"=52.25235201667956; "G"=3.731226274323352; "J"=10.38882643771372; "L"=11.59624026431456; "N"=21.30283811754617 }
    23 = @{ "B"=30.4896431837013; "C"=13.63455743655829; "E"=10.07506156052845; "F"=52.12823386169019; "G"=3.733113893424147; "J"=10.39926995030477; "L"=11.59483790188933; "N"=21.32366510638106 }
    24 = @{ "B"=29.7442159464253; "C"=12.81917355503657; "E"=10.03773529381916; "F"=51.6785410918277; "G"=3.740521075240168; "J"=10.4408567403073; "L"=11.59272890038292; "N"=21.40609471468834 }
    25 = @{ "B"=28.95237091175497; "C"=11.89134443317788; "E"=9.996928916878545; "F"=51.23831634227832; "G"=3.749068168753439; "J"=10.49005524251856; "L"=11.59733136582204; "N"=21.50256900031225 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
